$wb = $excel.ActiveWorkbook

# --- Sheet "data": append the 2024 annual row (row 99) ---
$wsData = $wb.Worksheets.Item("data")

# Copy the formatting of the last existing data row (98) down into the new row (99)
$wsData.Range("A98:G98").Copy()
$wsData.Range("A99:G99").PasteSpecial(-4122)

$wsData.Range("A99").Value = 2024
$wsData.Range("B99").Value = 0.19369429191554208
$wsData.Range("C99").Value = 0.22216568002960194
$wsData.Range("D99").Value = 0.008
$wsData.Range("E99").Value = 5881.63
$wsData.Range("F99").Value = 315.605
$wsData.Range("G99").Value = 74.83

# --- Sheet "readme": roll the "...1927-2023" labels forward to "...1927-2024" ---
$wsReadme = $wb.Worksheets.Item("readme")
$wsReadme.Range("A1").Value = "Annual Data 1927-2024"
$wsReadme.Range("B4").Value = "1928-2024"
$wsReadme.Range("B5").Value = "Final day year 1927-2024"
$wsReadme.Range("B6").Value = "December 1927-2024"
$wsReadme.Range("B7").Value = "1927-2024"

# --- Restore/update selection state on both sheets ---
$wsData.Range("D102").Select()
$wsReadme.Activate()
$wsReadme.Range("B8").Select()
